# Update result_data_KNN values (algorithm re-run produced slightly
# different imputed numbers for a subset of cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E3"  = 12.643
    "A4"  = -21.984
    "C4"  = -12.78
    "E4"  = 12.647
    "C5"  = -12.894
    "A6"  = -21.108
    "A7"  = -21.038
    "C8"  = -12.796
    "E9"  = 13.018
    "E11" = 12.852
    "E14" = 13.056
    "A16" = -21.061
    "C16" = -12.422
    "E18" = 12.657
    "A20" = -22.096
    "C22" = -12.78
    "E25" = 12.827
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
